$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new header cells I1 (I0) and J1 (IF)
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# Copy the header style (border/bold/alignment) from H1 onto I1:J1
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# Fill in the new I/J column values for each data row (rows 2-50)
$data = @(
    @(8, 8),
    @(6, 6),
    @(6, 6),
    @(7, 7),
    @(8, 8),
    @(7, 7),
    @(8, 9),
    @(8, 8),
    @(9, 9),
    @(9, 9),
    @(8, 8),
    @(9, 9),
    @(9, 9),
    @(9, 9),
    @(9, 9),
    @(8, 8),
    @(8, 8),
    @(8, 8),
    @(8, 8),
    @(7, 7),
    @(8, 8),
    @(7, 7),
    @(8, 8),
    @(5, 5),
    @(10, 10),
    @(7, 7),
    @(6, 6),
    @(5, 5),
    @(8, 8),
    @(8, 8),
    @(9, 9),
    @(6, 7),
    @(8, 8),
    @(8, 8),
    @(11, 12),
    @(9, 9),
    @(9, 9),
    @(3, 4),
    @(8, 8),
    @(7, 8),
    @(6, 6),
    @(5, 5),
    @(9, 9),
    @(7, 7),
    @(8, 8),
    @(4, 4),
    @(7, 7),
    @(5, 5),
    @(4, 4)
)

$row = 2
foreach ($pair in $data) {
    $ws.Cells.Item($row, 9).Value = $pair[0]
    $ws.Cells.Item($row, 10).Value = $pair[1]
    $row++
}
